# Add a second table of currency conversion rates below the existing
# header, processed/written into the worksheet (rows 3-8, columns A-D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$currencies = @(
    @("EUR", "euro",          1, "364,85"),
    @("BGN", "bolgár leva",   1, "186,54"),
    @("CZK", "cseh korona",   1, "13,35"),
    @("PLN", "lengyel zloty", 1, "80,02"),
    @("RON", "román lej",     1, "74,82"),
    @("UAH", "ukrán hrivnya", 1, "10,99")
)

$row = 3
foreach ($currency in $currencies) {
    $ws.Cells.Item($row, 1).Value = $currency[0]
    $ws.Cells.Item($row, 2).Value = $currency[1]
    $ws.Cells.Item($row, 3).Value = $currency[2]
    $ws.Cells.Item($row, 4).Value = $currency[3]
    $row++
}

$ws.Range("D8").Select() | Out-Null
